$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row's values first
$ws.Range("A16").Value = "K_CALCMETH"
$ws.Range("B16").Value = "Berechnungsmethode"
$ws.Range("C16").Value = "Calculation method"

# Copy the formatting (style) from the row above (row 15) onto the new
# row 16, mirroring the same cell style ("s=4") used by every other data
# row in the table, without disturbing the values we just set.
$ws.Range("A15:C15").Copy()
$ws.Range("A16:C16").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0
